$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 991.2308
$ws.Range("I18").Value = 534.9091
$ws.Range("J18").Value = 3501
$ws.Range("K18").Value = 534.9091
$ws.Range("L18").Value = 3501
$ws.Range("M18").Value = -250.9091
$ws.Range("N18").Value = -4069
$ws.Range("H28").Value = 565.0909
$ws.Range("I28").Value = 219.66667
$ws.Range("J28").Value = 1169.5834
$ws.Range("K28").Value = 219.66667
$ws.Range("L28").Value = 1169.5834
$ws.Range("M28").Value = 265.33333
$ws.Range("N28").Value = -2139.5834
$ws.Range("H58").Value = 714435.5
$ws.Range("I58").Value = 1639619.9
$ws.Range("J58").Value = 2755.1538
$ws.Range("K58").Value = 4918859.699999999
$ws.Range("L58").Value = 8265.4614
$ws.Range("M58").Value = -4918709.699999999
$ws.Range("N58").Value = -8565.4614
$ws.Range("H64").Value = 69612.336
$ws.Range("I64").Value = 252175
$ws.Range("J64").Value = 3225.9092
$ws.Range("K64").Value = 252175
$ws.Range("L64").Value = 3225.9092
$ws.Range("M64").Value = -251927
$ws.Range("N64").Value = -3721.9092
$ws.Range("H67").Value = 69612.336
$ws.Range("I67").Value = 252175
$ws.Range("J67").Value = 3225.9092
$ws.Range("K67").Value = 252175
$ws.Range("L67").Value = 3225.9092
$ws.Range("M67").Value = -251317
$ws.Range("N67").Value = -4941.9092
$ws.Range("H113").Value = 49409.523
$ws.Range("I113").Value = 101691
$ws.Range("J113").Value = 1880.909
$ws.Range("K113").Value = 101691
$ws.Range("L113").Value = 1880.909
$ws.Range("M113").Value = -98437
$ws.Range("N113").Value = -8388.909
$ws.Range("H116").Value = 4458.909
$ws.Range("I116").Value = 6599.6665
$ws.Range("J116").Value = 1890
$ws.Range("K116").Value = 6599.6665
$ws.Range("L116").Value = 1890
$ws.Range("M116").Value = -3157.6665
$ws.Range("N116").Value = -8774
$ws.Range("H132").Value = 8071732.5
$ws.Range("I132").Value = 9623367
$ws.Range("J132").Value = 3232.4
$ws.Range("K132").Value = 28870101
$ws.Range("L132").Value = 9697.200000000001
$ws.Range("M132").Value = -28867571
$ws.Range("N132").Value = -14757.2
$ws.Range("H137").Value = 2043.35
$ws.Range("I137").Value = 1391.6875
$ws.Range("K137").Value = 4175.0625
$ws.Range("M137").Value = -1625.0625
$ws.Range("H138").Value = 1744.7
$ws.Range("I138").Value = 543.86536
$ws.Range("J138").Value = 3045.6042
$ws.Range("K138").Value = 1631.59608
$ws.Range("L138").Value = 9136.812600000001
$ws.Range("M138").Value = 3508.40392
$ws.Range("N138").Value = -19416.8126
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1364.6765
$ws.Range("I74").Value = 696.7143
$ws.Range("J74").Value = 2443.6924
$ws.Range("K74").Value = 696.7143
$ws.Range("L74").Value = 2443.6924
$ws.Range("M74").Value = 177.2857
$ws.Range("N74").Value = -4191.6924
$ws.Range("H77").Value = 1364.6765
$ws.Range("I77").Value = 696.7143
$ws.Range("J77").Value = 2443.6924
$ws.Range("K77").Value = 3483.5715
$ws.Range("L77").Value = 12218.462
$ws.Range("M77").Value = 884.4285
$ws.Range("N77").Value = -20954.462
$ws.Range("H114").Value = 36295
$ws.Range("J114").Value = 36295
$ws.Range("L114").Value = 36295
$ws.Range("N114").Value = -44973
$ws.Range("H132").Value = 3902.1667
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3902.1667
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 11706.5001
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -16766.5001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 18599.6
$ws.Range("I26").Value = 18599.6
$ws.Range("K26").Value = 18599.6
$ws.Range("M26").Value = -18307.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1028.8334
$ws.Range("I22").Value = 334.57144
$ws.Range("J22").Value = 2000.8
$ws.Range("K22").Value = 334.57144
$ws.Range("L22").Value = 2000.8
$ws.Range("M22").Value = 15.42856
$ws.Range("N22").Value = -2700.8
$ws.Range("H31").Value = 26542.477
$ws.Range("I31").Value = 1007.8333
$ws.Range("J31").Value = 60588.668
$ws.Range("K31").Value = 1007.8333
$ws.Range("L31").Value = 60588.668
$ws.Range("M31").Value = -712.8333
$ws.Range("N31").Value = -61178.668
$ws.Range("H34").Value = 26542.477
$ws.Range("I34").Value = 1007.8333
$ws.Range("J34").Value = 60588.668
$ws.Range("K34").Value = 1007.8333
$ws.Range("L34").Value = 60588.668
$ws.Range("M34").Value = -805.8333
$ws.Range("N34").Value = -60992.668
$ws.Range("H35").Value = 2386.3635
$ws.Range("I35").Value = 792.8570999999999
$ws.Range("J35").Value = 5175
$ws.Range("K35").Value = 792.8570999999999
$ws.Range("L35").Value = 5175
$ws.Range("M35").Value = -498.8570999999999
$ws.Range("N35").Value = -5763
$ws.Range("H99").Value = 8558.294
$ws.Range("I99").Value = 3293.3333
$ws.Range("J99").Value = 9686.5
$ws.Range("K99").Value = 3293.3333
$ws.Range("L99").Value = 9686.5
$ws.Range("M99").Value = -1795.3333
$ws.Range("N99").Value = -12682.5
$ws.Range("H126").Value = 8558.294
$ws.Range("I126").Value = 3293.3333
$ws.Range("J126").Value = 9686.5
$ws.Range("K126").Value = 9879.999899999999
$ws.Range("L126").Value = 29059.5
$ws.Range("M126").Value = -7409.999899999999
$ws.Range("N126").Value = -33999.5
$ws.Range("H132").Value = 3237.2903
$ws.Range("I132").Value = 2962.1765
$ws.Range("J132").Value = 3571.3572
$ws.Range("K132").Value = 8886.529500000001
$ws.Range("L132").Value = 10714.0716
$ws.Range("M132").Value = -6356.529500000001
$ws.Range("N132").Value = -15774.0716
$ws.Range("H134").Value = 1317.2858
$ws.Range("I134").Value = 1221.0435
$ws.Range("K134").Value = 3663.1305
$ws.Range("M134").Value = -1128.1305
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 805.77
$ws.Range("J131").Value = 831.0319
$ws.Range("L131").Value = 2493.0957
$ws.Range("N131").Value = -12573.0957
$ws.Range("H132").Value = 1780.3636
$ws.Range("I132").Value = 950
$ws.Range("J132").Value = 2254.8572
$ws.Range("K132").Value = 8550
$ws.Range("L132").Value = 20293.7148
$ws.Range("M132").Value = -6020
$ws.Range("N132").Value = -25353.7148
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 34906.668
$ws.Range("J51").Value = 34906.668
$ws.Range("L51").Value = 34906.668
$ws.Range("N51").Value = -35924.668
$ws.Range("H102").Value = 1202496
$ws.Range("I102").Value = 3160
$ws.Range("K102").Value = 3160
$ws.Range("M102").Value = -1538
$ws.Range("H140").Value = 29526.666
$ws.Range("J140").Value = 29526.666
$ws.Range("L140").Value = 29526.666
$ws.Range("N140").Value = -39886.666
$ws.Range("H141").Value = 30645.8
$ws.Range("J141").Value = 30645.8
$ws.Range("L141").Value = 30645.8
$ws.Range("N141").Value = -41005.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2804.8604
$ws.Range("I132").Value = 2805.861
$ws.Range("J132").Value = 2799.7144
$ws.Range("K132").Value = 8417.582999999999
$ws.Range("L132").Value = 8399.143199999999
$ws.Range("M132").Value = -5887.582999999999
$ws.Range("N132").Value = -13459.1432
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 5838.4
$ws.Range("J41").Value = 5962.5
$ws.Range("L41").Value = 5962.5
$ws.Range("N41").Value = -6742.5
$ws.Range("H45").Value = 12498.333
$ws.Range("J45").Value = 12498.333
$ws.Range("L45").Value = 12498.333
$ws.Range("N45").Value = -13480.333
$ws.Range("H74").Value = 10309
$ws.Range("J74").Value = 10309
$ws.Range("L74").Value = 10309
$ws.Range("N74").Value = -12181
$ws.Range("H77").Value = 10309
$ws.Range("J77").Value = 10309
$ws.Range("L77").Value = 30927
$ws.Range("N77").Value = -40287
$ws.Range("H136").Value = 833
$ws.Range("I136").Value = 498.16666
$ws.Range("J136").Value = 1837.5
$ws.Range("K136").Value = 1494.49998
$ws.Range("L136").Value = 5512.5
$ws.Range("M136").Value = 1055.50002
$ws.Range("N136").Value = -10612.5
$ws.Range("H140").Value = 61250
$ws.Range("J140").Value = 61250
$ws.Range("L140").Value = 61250
$ws.Range("N140").Value = -71610
$ws.Range("H141").Value = 67707.5
$ws.Range("J141").Value = 67707.5
$ws.Range("L141").Value = 67707.5
$ws.Range("N141").Value = -78067.5
